$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns O (14) and P (15), matching N1 style/format
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

# Data rows 2-67: new O/P values
$ws.Range("O2").Value = -0.2230403794088271
$ws.Range("P2").Value = -0.2215874471047971
$ws.Range("O3").Value = 0.2371335201232264
$ws.Range("P3").Value = 0.2375663013777764
$ws.Range("O4").Value = 0.2253883573517043
$ws.Range("P4").Value = 0.2274855126263637
$ws.Range("O5").Value = -0.04374106398907703
$ws.Range("P5").Value = -0.04418221182924286
$ws.Range("O6").Value = 0.194987104306643
$ws.Range("P6").Value = 0.1950295963127417
$ws.Range("O7").Value = -0.3761258952878594
$ws.Range("P7").Value = -0.3745222279606223
$ws.Range("O8").Value = -0.1443539110529433
$ws.Range("P8").Value = -0.1409675007790592
$ws.Range("O9").Value = -0.32986226152445
$ws.Range("P9").Value = -0.3269832763153716
$ws.Range("O10").Value = 0.4065492170739562
$ws.Range("P10").Value = 0.4062837534899008
$ws.Range("O11").Value = -0.1767839794601572
$ws.Range("P11").Value = -0.1758867005252963
$ws.Range("O12").Value = -0.05340640561671622
$ws.Range("P12").Value = -0.05576876759396782
$ws.Range("O13").Value = -0.02267557700730256
$ws.Range("P13").Value = -0.02382135445782242
$ws.Range("O14").Value = 0.0869685881776623
$ws.Range("P14").Value = 0.08277770999501546
$ws.Range("O15").Value = -0.02450882852917694
$ws.Range("P15").Value = -0.03191811049595003
$ws.Range("O16").Value = 0.3720328181196113
$ws.Range("P16").Value = 0.3644626639804314
$ws.Range("O17").Value = 0.5051885710101555
$ws.Range("P17").Value = 0.4953917611788026
$ws.Range("O18").Value = -0.113170310948517
$ws.Range("P18").Value = -0.1183636471376141
$ws.Range("O19").Value = 0.3569138777307417
$ws.Range("P19").Value = 0.3519627866354487
$ws.Range("O20").Value = 0.2850369672763723
$ws.Range("P20").Value = 0.2760019920768016
$ws.Range("O21").Value = 0.5625723715019986
$ws.Range("P21").Value = 0.5552931054562932
$ws.Range("O22").Value = 0.3211604717349571
$ws.Range("P22").Value = 0.3129834570441234
$ws.Range("O23").Value = -0.09623675534332604
$ws.Range("P23").Value = -0.1033090105301266
$ws.Range("O24").Value = 2.515881345646089
$ws.Range("P24").Value = 2.385979903177596
$ws.Range("O25").Value = 0.3163499125572618
$ws.Range("P25").Value = 0.3142420485315187
$ws.Range("O26").Value = 0.2078480175611874
$ws.Range("P26").Value = 0.2008989692326813
$ws.Range("O27").Value = 0.07679443510443131
$ws.Range("P27").Value = 0.07116955244656895
$ws.Range("O28").Value = 0.8187656121063521
$ws.Range("P28").Value = 0.8157019941728602
$ws.Range("O29").Value = 2.205293729993913
$ws.Range("P29").Value = 2.109496983584984
$ws.Range("O30").Value = 0.696298317405343
$ws.Range("P30").Value = 0.690487154084958
$ws.Range("O31").Value = -0.4271221476294025
$ws.Range("P31").Value = -0.4275181353616064
$ws.Range("O32").Value = 0.5794123627511983
$ws.Range("P32").Value = 0.5753999204112504
$ws.Range("O33").Value = 0.7698048138358585
$ws.Range("P33").Value = 0.7678621658740948
$ws.Range("O34").Value = -0.7972732831639558
$ws.Range("P34").Value = -0.7992227449451159
$ws.Range("O35").Value = 0.8100912676638644
$ws.Range("P35").Value = 0.8109367755268172
$ws.Range("O36").Value = 0.7663594715324018
$ws.Range("P36").Value = 0.7686663851815934
$ws.Range("O37").Value = 0.7315327099919052
$ws.Range("P37").Value = 0.7337836788962876
$ws.Range("O38").Value = 0.6527207191898265
$ws.Range("P38").Value = 0.6524288853530678
$ws.Range("O39").Value = 0.6076357902385338
$ws.Range("P39").Value = 0.6090920531484927
$ws.Range("O40").Value = 0.7750233809010156
$ws.Range("P40").Value = 0.7760218781801852
$ws.Range("O41").Value = 0.5601800806008893
$ws.Range("P41").Value = 0.5619234513548153
$ws.Range("O42").Value = 0.6366265908810105
$ws.Range("P42").Value = 0.6377126728415672
$ws.Range("O43").Value = 0.6770268791486883
$ws.Range("P43").Value = 0.6776495845415114
$ws.Range("O44").Value = 0.6945023137754879
$ws.Range("P44").Value = 0.6965753555995946
$ws.Range("O45").Value = 0.6971387595924574
$ws.Range("P45").Value = 0.701335022599936
$ws.Range("O46").Value = -1.238250268119329
$ws.Range("P46").Value = -1.238747748000431
$ws.Range("O47").Value = -0.9607655640437177
$ws.Range("P47").Value = -0.9610963028519899
$ws.Range("O48").Value = -0.8315234676308305
$ws.Range("P48").Value = -0.8304727980913619
$ws.Range("O49").Value = -0.6036330839074165
$ws.Range("P49").Value = -0.6028902152372047
$ws.Range("O50").Value = -0.03934679208788293
$ws.Range("P50").Value = -0.03985078781523613
$ws.Range("O51").Value = -0.8166511562403003
$ws.Range("P51").Value = -0.8153143551828111
$ws.Range("O52").Value = -0.8166511562403003
$ws.Range("P52").Value = -0.8153143551828111
$ws.Range("O53").Value = -1.085601938429782
$ws.Range("P53").Value = -1.085126782784815
$ws.Range("O54").Value = -0.148659336199851
$ws.Range("P54").Value = -0.1482352769675265
$ws.Range("O55").Value = -0.9778533050628286
$ws.Range("P55").Value = -0.9783622487125703
$ws.Range("O56").Value = -0.8543084915818779
$ws.Range("P56").Value = -0.8555862503769061
$ws.Range("O57").Value = -0.9083913972674763
$ws.Range("P57").Value = -0.9108654501463399
$ws.Range("O58").Value = -1.055668530183284
$ws.Range("P58").Value = -1.057403945582217
$ws.Range("O59").Value = -0.7627363318599898
$ws.Range("P59").Value = -0.762100072594416
$ws.Range("O60").Value = -0.413079003709617
$ws.Range("P60").Value = -0.4132176533454865
$ws.Range("O61").Value = 0.3899920230016475
$ws.Range("P61").Value = 0.3901378341989215
$ws.Range("O62").Value = -1.158370936788779
$ws.Range("P62").Value = -1.161126665179075
$ws.Range("O63").Value = -0.5613799767931202
$ws.Range("P63").Value = -0.5577499429162007
$ws.Range("O64").Value = -0.8377181870774687
$ws.Range("P64").Value = -0.836640146480697
$ws.Range("O65").Value = -0.02655587603658779
$ws.Range("P65").Value = -0.02650718150294729
$ws.Range("O66").Value = -0.7407928953198633
$ws.Range("P66").Value = -0.7438824701942061
$ws.Range("O67").Value = -0.7053170509057659
$ws.Range("P67").Value = -0.7099591814859051
